# export-structure-stats.xlsx: add two new "orientation" rows to the
# "Répartition des orientations" block on the stats sheet
#   - row 102 (previously blank): "Orientation vers CIAS"
#   - a brand-new row 103: "Autre orientation"
# Inserting the row naturally pushes the rest of the block (the section
# header + the 10 interaction rows that follow it) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 103 — everything from row 103
# down (the "Répartition des orientations" header and its rows) shifts
# down by one, to rows 105-115.
$ws.Rows.Item(103).Insert()

# Former row 102 (B102) was an empty placeholder row; give it a label.
$ws.Cells.Item(102, 2).Value = "Orientation vers CIAS"
$ws.Rows.Item(102).RowHeight = 16

# The freshly-inserted row 103 gets the second new label.
$ws.Cells.Item(103, 2).Value = "Autre orientation"
$ws.Rows.Item(103).RowHeight = 16

# Reflect the author's new scroll position / selection in the sheet view.
$excel.ActiveWindow.ScrollRow = 72
$ws.Range("B101").Select()
